# Auto-generated cell value updates (scheduled market-data refresh)
$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 9316393
$ws.Range("I70").Value = 20960246
$ws.Range("J70").Value = 1310
$ws.Range("K70").Value = 62880738
$ws.Range("L70").Value = 3930
$ws.Range("M70").Value = -62880468
$ws.Range("N70").Value = -4470
$ws.Range("H73").Value = 9316393
$ws.Range("I73").Value = 20960246
$ws.Range("J73").Value = 1310
$ws.Range("K73").Value = 62880738
$ws.Range("L73").Value = 3930
$ws.Range("M73").Value = -62879802
$ws.Range("N73").Value = -5802
$ws.Range("H137").Value = 1077.875
$ws.Range("I137").Value = 884.5273
$ws.Range("J137").Value = 1703.4117
$ws.Range("K137").Value = 2653.5819
$ws.Range("L137").Value = 5110.2351
$ws.Range("M137").Value = -103.5819000000001
$ws.Range("N137").Value = -10210.2351
$ws.Range("H138").Value = 2621.8289
$ws.Range("I138").Value = 1671.3208
$ws.Range("J138").Value = 4812.1304
$ws.Range("K138").Value = 5013.9624
$ws.Range("L138").Value = 14436.3912
$ws.Range("M138").Value = 126.0375999999997
$ws.Range("N138").Value = -24716.3912
$ws.Range("H139").Value = 51500
$ws.Range("J139").Value = 51500
$ws.Range("L139").Value = 51500
$ws.Range("N139").Value = -61780
$ws.Range("H141").Value = 3723.9424
$ws.Range("I141").Value = 1728.5106
$ws.Range("J141").Value = 22481
$ws.Range("K141").Value = 5185.531800000001
$ws.Range("L141").Value = 67443
$ws.Range("M141").Value = -5.531800000000658
$ws.Range("N141").Value = -77803

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 957.9648999999999
$ws.Range("I61").Value = 935.4894
$ws.Range("K61").Value = 935.4894
$ws.Range("M61").Value = -723.4894
$ws.Range("H74").Value = 991.7954999999999
$ws.Range("I74").Value = 1023.5946
$ws.Range("J74").Value = 823.7143
$ws.Range("K74").Value = 1023.5946
$ws.Range("L74").Value = 823.7143
$ws.Range("M74").Value = -149.5946
$ws.Range("N74").Value = -2571.7143
$ws.Range("H77").Value = 991.7954999999999
$ws.Range("I77").Value = 1023.5946
$ws.Range("J77").Value = 823.7143
$ws.Range("K77").Value = 5117.973
$ws.Range("L77").Value = 4118.5715
$ws.Range("M77").Value = -749.973
$ws.Range("N77").Value = -12854.5715
$ws.Range("H136").Value = 957.9648999999999
$ws.Range("I136").Value = 935.4894
$ws.Range("K136").Value = 2806.4682
$ws.Range("M136").Value = -256.4682000000003

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 2007.5918
$ws.Range("I134").Value = 1540.742
$ws.Range("J134").Value = 2811.611
$ws.Range("K134").Value = 4622.226
$ws.Range("L134").Value = 8434.832999999999
$ws.Range("M134").Value = -2087.226
$ws.Range("N134").Value = -13504.833

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 266356.53
$ws.Range("I132").Value = 347663.66
$ws.Range("J132").Value = 2108.4167
$ws.Range("K132").Value = 1042990.98
$ws.Range("L132").Value = 6325.250100000001
$ws.Range("M132").Value = -1040460.98
$ws.Range("N132").Value = -11385.2501

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 771.375
$ws.Range("I2").Value = 1772.8334
$ws.Range("J2").Value = 170.5
$ws.Range("K2").Value = 10637.0004
$ws.Range("L2").Value = 1023
$ws.Range("M2").Value = -10524.0004
$ws.Range("N2").Value = -1249
$ws.Range("H11").Value = 247.71428
$ws.Range("I11").Value = 83.5
$ws.Range("J11").Value = 466.66666
$ws.Range("K11").Value = 250.5
$ws.Range("L11").Value = 1399.99998
$ws.Range("M11").Value = -110.5
$ws.Range("N11").Value = -1679.99998
$ws.Range("H113").Value = 549.75
$ws.Range("I113").Value = 484.5625
$ws.Range("J113").Value = 636.6667
$ws.Range("K113").Value = 1453.6875
$ws.Range("L113").Value = 1910.0001
$ws.Range("M113").Value = 716.3125
$ws.Range("N113").Value = -6250.0001
$ws.Range("H131").Value = 10528918
$ws.Range("J131").Value = 10754163
$ws.Range("L131").Value = 32262489
$ws.Range("N131").Value = -32272569

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 7005.8237
$ws.Range("I70").Value = 6255.5557
$ws.Range("J70").Value = 7849.875
$ws.Range("K70").Value = 6255.5557
$ws.Range("L70").Value = 7849.875
$ws.Range("M70").Value = -5985.5557
$ws.Range("N70").Value = -8389.875
$ws.Range("H73").Value = 7005.8237
$ws.Range("I73").Value = 6255.5557
$ws.Range("J73").Value = 7849.875
$ws.Range("K73").Value = 6255.5557
$ws.Range("L73").Value = 7849.875
$ws.Range("M73").Value = -5319.5557
$ws.Range("N73").Value = -9721.875
$ws.Range("H95").Value = 1125965.4
$ws.Range("J95").Value = 1125965.4
$ws.Range("L95").Value = 1125965.4
$ws.Range("N95").Value = -1131457.4
$ws.Range("H132").Value = 1822.6666
$ws.Range("I132").Value = 1214
$ws.Range("J132").Value = 3496.5
$ws.Range("K132").Value = 3642
$ws.Range("L132").Value = 10489.5
$ws.Range("M132").Value = -1112
$ws.Range("N132").Value = -15549.5
$ws.Range("H140").Value = 90369.75
$ws.Range("I140").Value = 60709
$ws.Range("J140").Value = 100256.664
$ws.Range("K140").Value = 60709
$ws.Range("L140").Value = 100256.664
$ws.Range("M140").Value = -55529
$ws.Range("N140").Value = -110616.664
$ws.Range("H141").Value = 70214.5
$ws.Range("J141").Value = 70214.5
$ws.Range("L141").Value = 70214.5
$ws.Range("N141").Value = -80574.5

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2022.3158
$ws.Range("I132").Value = 1911.6586
$ws.Range("J132").Value = 2305.875
$ws.Range("K132").Value = 5734.9758
$ws.Range("L132").Value = 6917.625
$ws.Range("M132").Value = -3204.9758
$ws.Range("N132").Value = -11977.625
$ws.Range("H136").Value = 1489.0521
$ws.Range("I136").Value = 1318.8334
$ws.Range("J136").Value = 2226.6667
$ws.Range("K136").Value = 3956.5002
$ws.Range("L136").Value = 6680.000100000001
$ws.Range("M136").Value = -1406.5002
$ws.Range("N136").Value = -11780.0001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5001.5
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 5001.5
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 5001.5
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -6249.5
$ws.Range("H65").Value = 5001.5
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 5001.5
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 25007.5
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -31247.5
$ws.Range("H75").Value = 56104
$ws.Range("I75").Value = 200000
$ws.Range("J75").Value = 20130
$ws.Range("K75").Value = 200000
$ws.Range("L75").Value = 20130
$ws.Range("M75").Value = -199064
$ws.Range("N75").Value = -22002
$ws.Range("H78").Value = 56104
$ws.Range("I78").Value = 200000
$ws.Range("J78").Value = 20130
$ws.Range("K78").Value = 600000
$ws.Range("L78").Value = 60390
$ws.Range("M78").Value = -595320
$ws.Range("N78").Value = -69750
$ws.Range("H86").Value = 134260
$ws.Range("J86").Value = 134260
$ws.Range("L86").Value = 134260
$ws.Range("N86").Value = -136506
$ws.Range("H89").Value = 134260
$ws.Range("J89").Value = 134260
$ws.Range("L89").Value = 671300
$ws.Range("N89").Value = -682532
$ws.Range("H92").Value = 31666.666
$ws.Range("J92").Value = 31666.666
$ws.Range("L92").Value = 31666.666
$ws.Range("N92").Value = -36658.666
$ws.Range("H132").Value = 786.4308
$ws.Range("I132").Value = 597.1
$ws.Range("J132").Value = 1417.5333
$ws.Range("K132").Value = 1791.3
$ws.Range("L132").Value = 4252.5999
$ws.Range("M132").Value = 738.6999999999998
$ws.Range("N132").Value = -9312.599900000001
